# Deploying to gh-pages — add 2019/2020 data columns (P, Q) to the table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (year headers) : copy formatting from the existing year header cells ---
$ws.Range("O4").Copy()
$ws.Range("P4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("P4").Value = 2019

$ws.Range("O4").Copy()
$ws.Range("Q4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("Q4").Value = 2020

# --- Row 5 (share of renewables, percent) : use the "0.0" formatted cell style ---
$ws.Range("H5").Copy()
$ws.Range("P5").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("P5").Value = 35.67

$ws.Range("H5").Copy()
$ws.Range("Q5").PasteSpecial(-4122)   # xlPasteFormats, value intentionally left blank

# --- Row 6 (hydropower electricity production, mln kWh) ---
$ws.Range("O6").Copy()
$ws.Range("P6").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("P6").Value = 13859.3

$ws.Range("O6").Copy()
$ws.Range("Q6").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("Q6").Value = 13979.1

# Clear clipboard/marching-ants selection state left over from copy operations.
$excel.CutCopyMode = $false

# Mirror the author's final cursor position recorded in the saved file.
$ws.Range("P9").Select()
